{"js": "// The document body contains a single 5-column table of two-digit\n// division problems. Every 4th row (0, 4, 8, 12, 16) holds the actual\n// \"a\u00f7b=c, d\" answer text; the rows between are spacer rows. Several of\n// the \"before\" strings repeat across different cells (e.g. \"78\u00f73=26, 0\"\n// and \"69\u00f75=13, 4\" each occur twice), so a global text search/replace\n// would be ambiguous. Instead we target each cell by its (row, column)\n// position and overwrite only the run text inside it, which leaves the\n// existing paragraph/run formatting (font, size, alignment) untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"91\u00f77=13, 0\" },\n  { row: 0, col: 1, text: \"68\u00f72=34, 0\" },\n  { row: 0, col: 2, text: \"59\u00f72=29, 1\" },\n  { row: 0, col: 3, text: \"17\u00f79=1, 8\" },\n  { row: 0, col: 4, text: \"24\u00f78=3, 0\" },\n  { row: 4, col: 0, text: \"79\u00f74=19, 3\" },\n  { row: 4, col: 1, text: \"80\u00f78=10, 0\" },\n  { row: 4, col: 2, text: \"16\u00f74=4, 0\" },\n  { row: 4, col: 3, text: \"37\u00f79=4, 1\" },\n  { row: 4, col: 4, text: \"34\u00f78=4, 2\" },\n  { row: 8, col: 0, text: \"61\u00f73=20, 1\" },\n  { row: 8, col: 1, text: \"61\u00f77=8, 5\" },\n  { row: 8, col: 2, text: \"93\u00f72=46, 1\" },\n  { row: 8, col: 3, text: \"41\u00f72=20, 1\" },\n  { row: 8, col: 4, text: \"86\u00f77=12, 2\" },\n  { row: 12, col: 0, text: \"87\u00f75=17, 2\" },\n  { row: 12, col: 1, text: \"90\u00f78=11, 2\" },\n  { row: 12, col: 2, text: \"21\u00f75=4, 1\" },\n  { row: 12, col: 3, text: \"62\u00f74=15, 2\" },\n  { row: 12, col: 4, text: \"18\u00f75=3, 3\" },\n  { row: 16, col: 0, text: \"54\u00f72=27, 0\" },\n  { row: 16, col: 1, text: \"83\u00f79=9, 2\" },\n  { row: 16, col: 2, text: \"71\u00f73=23, 2\" },\n  { row: 16, col: 3, text: \"70\u00f73=23, 1\" },\n  { row: 16, col: 4, text: \"89\u00f77=12, 5\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  // Replace just the text of the cell's (only) paragraph so the run's\n  // rPr (font/size) and the paragraph's pPr (alignment) survive as-is.\n  const para = cell.body.paragraphs.getFirst();\n  const range = para.getRange();\n  range.insertText(u.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "# The document body contains a single 5-column table of two-digit\n# division problems. Every 4th row (rows 1, 5, 9, 13, 17 in 1-based COM\n# indexing) holds the actual \"a\u00f7b=c, d\" answer text; the rows between\n# are spacer rows. Several \"before\" strings repeat across different\n# cells (e.g. \"78\u00f73=26, 0\" and \"69\u00f75=13, 4\" each occur twice), so a\n# global Find/Replace would be ambiguous. Instead we target each cell\n# directly by its (row, column) position via Table.Cell(row, column)\n# and overwrite only Range.Text, which keeps the existing run/paragraph\n# formatting (font, size, alignment) intact.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$table.Cell(1, 1).Range.Text = \"91\u00f77=13, 0\"\n$table.Cell(1, 2).Range.Text = \"68\u00f72=34, 0\"\n$table.Cell(1, 3).Range.Text = \"59\u00f72=29, 1\"\n$table.Cell(1, 4).Range.Text = \"17\u00f79=1, 8\"\n$table.Cell(1, 5).Range.Text = \"24\u00f78=3, 0\"\n$table.Cell(5, 1).Range.Text = \"79\u00f74=19, 3\"\n$table.Cell(5, 2).Range.Text = \"80\u00f78=10, 0\"\n$table.Cell(5, 3).Range.Text = \"16\u00f74=4, 0\"\n$table.Cell(5, 4).Range.Text = \"37\u00f79=4, 1\"\n$table.Cell(5, 5).Range.Text = \"34\u00f78=4, 2\"\n$table.Cell(9, 1).Range.Text = \"61\u00f73=20, 1\"\n$table.Cell(9, 2).Range.Text = \"61\u00f77=8, 5\"\n$table.Cell(9, 3).Range.Text = \"93\u00f72=46, 1\"\n$table.Cell(9, 4).Range.Text = \"41\u00f72=20, 1\"\n$table.Cell(9, 5).Range.Text = \"86\u00f77=12, 2\"\n$table.Cell(13, 1).Range.Text = \"87\u00f75=17, 2\"\n$table.Cell(13, 2).Range.Text = \"90\u00f78=11, 2\"\n$table.Cell(13, 3).Range.Text = \"21\u00f75=4, 1\"\n$table.Cell(13, 4).Range.Text = \"62\u00f74=15, 2\"\n$table.Cell(13, 5).Range.Text = \"18\u00f75=3, 3\"\n$table.Cell(17, 1).Range.Text = \"54\u00f72=27, 0\"\n$table.Cell(17, 2).Range.Text = \"83\u00f79=9, 2\"\n$table.Cell(17, 3).Range.Text = \"71\u00f73=23, 2\"\n$table.Cell(17, 4).Range.Text = \"70\u00f73=23, 1\"\n$table.Cell(17, 5).Range.Text = \"89\u00f77=12, 5\"\n"}
